$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Training trials stimuli replaced: rows 2-5 (B:D) now all point at the
# same non-social stimuli set instead of the per-row Mandalas_new files.
$ws.Range("B2").Value = "Stimuli/neutral/Still.jpg"
$ws.Range("C2").Value = "Stimuli/nonsocial/Right Tick.mp4"
$ws.Range("D2").Value = "Stimuli/nonsocial/Wrong Tick.mp4"

$ws.Range("B3").Value = "Stimuli/neutral/Still.jpg"
$ws.Range("C3").Value = "Stimuli/nonsocial/Right Tick.mp4"
$ws.Range("D3").Value = "Stimuli/nonsocial/Wrong Tick.mp4"

$ws.Range("B4").Value = "Stimuli/neutral/Still.jpg"
$ws.Range("C4").Value = "Stimuli/nonsocial/Right Tick.mp4"
$ws.Range("D4").Value = "Stimuli/nonsocial/Wrong Tick.mp4"

$ws.Range("B5").Value = "Stimuli/neutral/Still.jpg"
$ws.Range("C5").Value = "Stimuli/nonsocial/Right Tick.mp4"
$ws.Range("D5").Value = "Stimuli/nonsocial/Wrong Tick.mp4"

# Selection moved from F10 to D6 in the saved view state.
$ws.Range("D6").Select()
